$wb = $excel.ActiveWorkbook

# --- Matriz_Resultados ---
$ws = $wb.Worksheets.Item("Matriz_Resultados")
$data = New-Object 'object[,]' 9,9
$data[0,0] = 0
$data[0,1] = -1
$data[0,2] = 0
$data[0,3] = 0
$data[0,4] = 0
$data[0,5] = 0
$data[0,6] = 0
$data[0,7] = 0
$data[0,8] = 0
$data[1,0] = 1
$data[1,1] = 0
$data[1,2] = 1
$data[1,3] = 0
$data[1,4] = 1
$data[1,5] = 1
$data[1,6] = 1
$data[1,7] = 1
$data[1,8] = 1
$data[2,0] = 0
$data[2,1] = -1
$data[2,2] = 0
$data[2,3] = 0
$data[2,4] = 0
$data[2,5] = 0
$data[2,6] = 0
$data[2,7] = -1
$data[2,8] = 1
$data[3,0] = 0
$data[3,1] = 0
$data[3,2] = 0
$data[3,3] = 0
$data[3,4] = 0
$data[3,5] = 0
$data[3,6] = 0
$data[3,7] = 0
$data[3,8] = 0
$data[4,0] = 0
$data[4,1] = -1
$data[4,2] = 0
$data[4,3] = 0
$data[4,4] = 0
$data[4,5] = 0
$data[4,6] = 0
$data[4,7] = -1
$data[4,8] = 0
$data[5,0] = 0
$data[5,1] = -1
$data[5,2] = 0
$data[5,3] = 0
$data[5,4] = 0
$data[5,5] = 0
$data[5,6] = 0
$data[5,7] = -1
$data[5,8] = 1
$data[6,0] = 0
$data[6,1] = -1
$data[6,2] = 0
$data[6,3] = 0
$data[6,4] = 0
$data[6,5] = 0
$data[6,6] = 0
$data[6,7] = -1
$data[6,8] = 1
$data[7,0] = 0
$data[7,1] = -1
$data[7,2] = 1
$data[7,3] = 0
$data[7,4] = 1
$data[7,5] = 1
$data[7,6] = 1
$data[7,7] = 0
$data[7,8] = 1
$data[8,0] = 0
$data[8,1] = -1
$data[8,2] = -1
$data[8,3] = 0
$data[8,4] = 0
$data[8,5] = -1
$data[8,6] = -1
$data[8,7] = -1
$data[8,8] = 0
$ws.Range("B2:J10").Value = $data

# --- P_valores ---
$ws = $wb.Worksheets.Item("P_valores")
$data = New-Object 'object[,]' 9,9
$data[0,0] = 1
$data[0,1] = 0.001245405197254135
$data[0,2] = 0.02240938001384807
$data[0,3] = 0.2434043997816189
$data[0,4] = 0.2524716231943038
$data[0,5] = 0.0172075437388346
$data[0,6] = 0.01004846552458516
$data[0,7] = 0.001840280335065669
$data[0,8] = 0.004942283565851469
$data[1,0] = 0.001245405197254135
$data[1,1] = 1
$data[1,2] = [double]"5.043989448338948E-09"
$data[1,3] = 0.001520976891660197
$data[1,4] = 0.0004757001858719523
$data[1,5] = [double]"1.391435199726843E-06"
$data[1,6] = [double]"5.87940147123156E-06"
$data[1,7] = 0.0006760435113633001
$data[1,8] = [double]"3.239860466131006E-06"
$data[2,0] = 0.02240938001384807
$data[2,1] = [double]"5.043989448338948E-09"
$data[2,2] = 1
$data[2,3] = 0.01424890445561555
$data[2,4] = 0.008849100021360412
$data[2,5] = 0.8004402842736695
$data[2,6] = 0.5823195383159661
$data[2,7] = [double]"7.18838033542113E-09"
$data[2,8] = [double]"4.538992476743609E-05"
$data[3,0] = 0.2434043997816189
$data[3,1] = 0.001520976891660197
$data[3,2] = 0.01424890445561555
$data[3,3] = 1
$data[3,4] = 0.3638407383778444
$data[3,5] = 0.01755960468269802
$data[3,6] = 0.01147569869285303
$data[3,7] = 0.002049300521915765
$data[3,8] = 0.5249190524081828
$data[4,0] = 0.2524716231943038
$data[4,1] = 0.0004757001858719523
$data[4,2] = 0.008849100021360412
$data[4,3] = 0.3638407383778444
$data[4,4] = 1
$data[4,5] = 0.006642593734019098
$data[4,6] = 0.003779803164086104
$data[4,7] = 0.0007265070874398205
$data[4,8] = 0.00930434681100123
$data[5,0] = 0.0172075437388346
$data[5,1] = [double]"1.391435199726843E-06"
$data[5,2] = 0.8004402842736695
$data[5,3] = 0.01755960468269802
$data[5,4] = 0.006642593734019098
$data[5,5] = 1
$data[5,6] = 0.3195704027182673
$data[5,7] = [double]"5.755291052178535E-06"
$data[5,8] = [double]"3.908529823970497E-05"
$data[6,0] = 0.01004846552458516
$data[6,1] = [double]"5.87940147123156E-06"
$data[6,2] = 0.5823195383159661
$data[6,3] = 0.01147569869285303
$data[6,4] = 0.003779803164086104
$data[6,5] = 0.3195704027182673
$data[6,6] = 1
$data[6,7] = [double]"1.105020764158304E-05"
$data[6,8] = [double]"2.256782427201287E-05"
$data[7,0] = 0.001840280335065669
$data[7,1] = 0.0006760435113633001
$data[7,2] = [double]"7.18838033542113E-09"
$data[7,3] = 0.002049300521915765
$data[7,4] = 0.0007265070874398205
$data[7,5] = [double]"5.755291052178535E-06"
$data[7,6] = [double]"1.105020764158304E-05"
$data[7,7] = 1
$data[7,8] = [double]"4.677304772382485E-06"
$data[8,0] = 0.004942283565851469
$data[8,1] = [double]"3.239860466131006E-06"
$data[8,2] = [double]"4.538992476743609E-05"
$data[8,3] = 0.5249190524081828
$data[8,4] = 0.00930434681100123
$data[8,5] = [double]"3.908529823970497E-05"
$data[8,6] = [double]"2.256782427201287E-05"
$data[8,7] = [double]"4.677304772382485E-06"
$data[8,8] = 1
$ws.Range("B2:J10").Value = $data

# --- Estadisticos_DM ---
$ws = $wb.Worksheets.Item("Estadisticos_DM")
$data = New-Object 'object[,]' 9,9
$data[0,0] = 0
$data[0,1] = 4.02815138357157
$data[0,2] = 2.56602960057472
$data[0,3] = -1.21786281698874
$data[0,4] = -1.193592427711046
$data[0,5] = 2.701455150592083
$data[0,6] = 2.974401866180414
$data[0,7] = 3.829552381439773
$data[0,8] = -3.331531058401805
$data[1,0] = -4.02815138357157
$data[1,1] = 0
$data[1,2] = -12.58908705870013
$data[1,3] = -3.926303524301437
$data[1,4] = -4.525203327505016
$data[1,5] = -7.989063045642583
$data[1,6] = -7.038480500696632
$data[1,7] = -4.342239038476555
$data[1,8] = -7.422396709428656
$data[2,0] = -2.56602960057472
$data[2,1] = 12.58908705870013
$data[2,2] = 0
$data[2,3] = -2.797548512407319
$data[2,4] = -3.038519125592238
$data[2,5] = -0.2576305742763857
$data[2,6] = 0.5630303453468123
$data[2,7] = 12.24823667788488
$data[2,8] = -5.808047212796896
$data[3,0] = 1.21786281698874
$data[3,1] = 3.926303524301437
$data[3,2] = 2.797548512407319
$data[3,3] = 0
$data[3,4] = 0.93861486008703
$data[3,5] = 2.691110712780848
$data[3,6] = 2.907282146497533
$data[3,7] = 3.775064298219911
$data[3,8] = -0.6520679940372778
$data[4,0] = 1.193592427711046
$data[4,1] = 4.525203327505016
$data[4,2] = 3.038519125592238
$data[4,3] = -0.93861486008703
$data[4,4] = 0
$data[4,5] = 3.182902667460546
$data[4,6] = 3.466343396011214
$data[4,7] = 4.304982913625855
$data[4,8] = -3.013225260221992
$data[5,0] = -2.701455150592083
$data[5,1] = 7.989063045642583
$data[5,2] = 0.2576305742763857
$data[5,3] = -2.691110712780848
$data[5,4] = -3.182902667460546
$data[5,5] = 0
$data[5,6] = 1.032002597761542
$data[5,7] = 7.052013475001109
$data[5,8] = -5.893993492430864
$data[6,0] = -2.974401866180414
$data[6,1] = 7.038480500696632
$data[6,2] = -0.5630303453468123
$data[6,3] = -2.907282146497533
$data[6,4] = -3.466343396011214
$data[6,5] = -1.032002597761542
$data[6,6] = 0
$data[6,7] = 6.645042906309448
$data[6,8] = -6.21484374659826
$data[7,0] = -3.829552381439773
$data[7,1] = 4.342239038476555
$data[7,2] = -12.24823667788488
$data[7,3] = -3.775064298219911
$data[7,4] = -4.304982913625855
$data[7,5] = -7.052013475001109
$data[7,6] = -6.645042906309448
$data[7,7] = 0
$data[7,8] = -7.184374872745349
$data[8,0] = 3.331531058401805
$data[8,1] = 7.422396709428656
$data[8,2] = 5.808047212796896
$data[8,3] = 0.6520679940372778
$data[8,4] = 3.013225260221992
$data[8,5] = 5.893993492430864
$data[8,6] = 6.21484374659826
$data[8,7] = 7.184374872745349
$data[8,8] = 0
$ws.Range("B2:J10").Value = $data

# --- Resumen ---
$ws = $wb.Worksheets.Item("Resumen")
$data = New-Object 'object[,]' 9,5
$data[0,0] = 7
$data[0,1] = 0
$data[0,2] = 1
$data[0,3] = 87.5
$data[0,4] = 0.548241150447542
$data[1,0] = 5
$data[1,1] = 1
$data[1,2] = 2
$data[1,3] = 62.5
$data[1,4] = 0.578071245476774
$data[2,0] = 1
$data[2,1] = 2
$data[2,2] = 5
$data[2,3] = 12.5
$data[2,4] = 0.7043433626341646
$data[3,0] = 1
$data[3,1] = 2
$data[3,2] = 5
$data[3,3] = 12.5
$data[3,4] = 0.719684280702942
$data[4,0] = 1
$data[4,1] = 2
$data[4,2] = 5
$data[4,3] = 12.5
$data[4,4] = 0.7147028466257251
$data[5,0] = 0
$data[5,1] = 1
$data[5,2] = 7
$data[5,3] = 0
$data[5,4] = 0.9466535455977692
$data[6,0] = 0
$data[6,1] = 0
$data[6,2] = 8
$data[6,3] = 0
$data[6,4] = 1.054448770626301
$data[7,0] = 0
$data[7,1] = 2
$data[7,2] = 6
$data[7,3] = 0
$data[7,4] = 0.9701992056259899
$data[8,0] = 0
$data[8,1] = 5
$data[8,2] = 3
$data[8,3] = 0
$data[8,4] = 1.105098924108501
$ws.Range("B2:F10").Value = $data

$ws.Range("A4").Value = "AV-MCPS"
$ws.Range("A5").Value = "MCPS"
